# Mise a jour avec les donnees categorielle.
# Columns N (Rentabilite), O (Taux Cout des produits vendus), P (Taux Opex)
# and Q (Taux Marge brute) for data rows 2-22 are switched from rounded
# numeric fractions (e.g. 0.5) to precise percentage text labels
# (e.g. "52.3%"), stored as text/categorical values rather than numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column N, rows 2-22
$cell = $ws.Range("N2")
$cell.NumberFormat = "@"
$cell.Value = "52.3%"
$cell.Style = "Normal"
$cell = $ws.Range("N3")
$cell.NumberFormat = "@"
$cell.Value = "5.2%"
$cell.Style = "Normal"
$cell = $ws.Range("N4")
$cell.NumberFormat = "@"
$cell.Value = "62.6%"
$cell.Style = "Normal"
$cell = $ws.Range("N5")
$cell.NumberFormat = "@"
$cell.Value = "27.0%"
$cell.Style = "Normal"
$cell = $ws.Range("N6")
$cell.NumberFormat = "@"
$cell.Value = "31.8%"
$cell.Style = "Normal"
$cell = $ws.Range("N7")
$cell.NumberFormat = "@"
$cell.Value = "-10.5%"
$cell.Style = "Normal"
$cell = $ws.Range("N8")
$cell.NumberFormat = "@"
$cell.Value = "42.7%"
$cell.Style = "Normal"
$cell = $ws.Range("N9")
$cell.NumberFormat = "@"
$cell.Value = "5.0%"
$cell.Style = "Normal"
$cell = $ws.Range("N10")
$cell.NumberFormat = "@"
$cell.Value = "62.0%"
$cell.Style = "Normal"
$cell = $ws.Range("N11")
$cell.NumberFormat = "@"
$cell.Value = "65.6%"
$cell.Style = "Normal"
$cell = $ws.Range("N12")
$cell.NumberFormat = "@"
$cell.Value = "10.8%"
$cell.Style = "Normal"
$cell = $ws.Range("N13")
$cell.NumberFormat = "@"
$cell.Value = "59.6%"
$cell.Style = "Normal"
$cell = $ws.Range("N14")
$cell.NumberFormat = "@"
$cell.Value = "59.5%"
$cell.Style = "Normal"
$cell = $ws.Range("N15")
$cell.NumberFormat = "@"
$cell.Value = "3.7%"
$cell.Style = "Normal"
$cell = $ws.Range("N16")
$cell.NumberFormat = "@"
$cell.Value = "53.0%"
$cell.Style = "Normal"
$cell = $ws.Range("N17")
$cell.NumberFormat = "@"
$cell.Value = "43.7%"
$cell.Style = "Normal"
$cell = $ws.Range("N18")
$cell.NumberFormat = "@"
$cell.Value = "17.7%"
$cell.Style = "Normal"
$cell = $ws.Range("N19")
$cell.NumberFormat = "@"
$cell.Value = "69.8%"
$cell.Style = "Normal"
$cell = $ws.Range("N20")
$cell.NumberFormat = "@"
$cell.Value = "64.7%"
$cell.Style = "Normal"
$cell = $ws.Range("N21")
$cell.NumberFormat = "@"
$cell.Value = "31.3%"
$cell.Style = "Normal"
$cell = $ws.Range("N22")
$cell.NumberFormat = "@"
$cell.Value = "76.5%"
$cell.Style = "Normal"

# Column O, rows 2-22
$cell = $ws.Range("O2")
$cell.NumberFormat = "@"
$cell.Value = "42.4%"
$cell.Style = "Normal"
$cell = $ws.Range("O3")
$cell.NumberFormat = "@"
$cell.Value = "63.5%"
$cell.Style = "Normal"
$cell = $ws.Range("O4")
$cell.NumberFormat = "@"
$cell.Value = "18.6%"
$cell.Style = "Normal"
$cell = $ws.Range("O5")
$cell.NumberFormat = "@"
$cell.Value = "64.0%"
$cell.Style = "Normal"
$cell = $ws.Range("O6")
$cell.NumberFormat = "@"
$cell.Value = "42.3%"
$cell.Style = "Normal"
$cell = $ws.Range("O7")
$cell.NumberFormat = "@"
$cell.Value = "56.1%"
$cell.Style = "Normal"
$cell = $ws.Range("O8")
$cell.NumberFormat = "@"
$cell.Value = "51.1%"
$cell.Style = "Normal"
$cell = $ws.Range("O9")
$cell.NumberFormat = "@"
$cell.Value = "64.9%"
$cell.Style = "Normal"
$cell = $ws.Range("O10")
$cell.NumberFormat = "@"
$cell.Value = "20.0%"
$cell.Style = "Normal"
$cell = $ws.Range("O11")
$cell.NumberFormat = "@"
$cell.Value = "27.5%"
$cell.Style = "Normal"
$cell = $ws.Range("O12")
$cell.NumberFormat = "@"
$cell.Value = "58.2%"
$cell.Style = "Normal"
$cell = $ws.Range("O13")
$cell.NumberFormat = "@"
$cell.Value = "15.4%"
$cell.Style = "Normal"
$cell = $ws.Range("O14")
$cell.NumberFormat = "@"
$cell.Value = "28.9%"
$cell.Style = "Normal"
$cell = $ws.Range("O15")
$cell.NumberFormat = "@"
$cell.Value = "63.7%"
$cell.Style = "Normal"
$cell = $ws.Range("O16")
$cell.NumberFormat = "@"
$cell.Value = "16.3%"
$cell.Style = "Normal"
$cell = $ws.Range("O17")
$cell.NumberFormat = "@"
$cell.Value = "50.4%"
$cell.Style = "Normal"
$cell = $ws.Range("O18")
$cell.NumberFormat = "@"
$cell.Value = "54.4%"
$cell.Style = "Normal"
$cell = $ws.Range("O19")
$cell.NumberFormat = "@"
$cell.Value = "15.4%"
$cell.Style = "Normal"
$cell = $ws.Range("O20")
$cell.NumberFormat = "@"
$cell.Value = "31.6%"
$cell.Style = "Normal"
$cell = $ws.Range("O21")
$cell.NumberFormat = "@"
$cell.Value = "45.7%"
$cell.Style = "Normal"
$cell = $ws.Range("O22")
$cell.NumberFormat = "@"
$cell.Value = "11.8%"
$cell.Style = "Normal"

# Column P, rows 2-22
$cell = $ws.Range("P2")
$cell.NumberFormat = "@"
$cell.Value = "5.3%"
$cell.Style = "Normal"
$cell = $ws.Range("P3")
$cell.NumberFormat = "@"
$cell.Value = "31.3%"
$cell.Style = "Normal"
$cell = $ws.Range("P4")
$cell.NumberFormat = "@"
$cell.Value = "18.8%"
$cell.Style = "Normal"
$cell = $ws.Range("P5")
$cell.NumberFormat = "@"
$cell.Value = "9.0%"
$cell.Style = "Normal"
$cell = $ws.Range("P6")
$cell.NumberFormat = "@"
$cell.Value = "25.9%"
$cell.Style = "Normal"
$cell = $ws.Range("P7")
$cell.NumberFormat = "@"
$cell.Value = "54.4%"
$cell.Style = "Normal"
$cell = $ws.Range("P8")
$cell.NumberFormat = "@"
$cell.Value = "6.2%"
$cell.Style = "Normal"
$cell = $ws.Range("P9")
$cell.NumberFormat = "@"
$cell.Value = "30.0%"
$cell.Style = "Normal"
$cell = $ws.Range("P10")
$cell.NumberFormat = "@"
$cell.Value = "18.0%"
$cell.Style = "Normal"
$cell = $ws.Range("P11")
$cell.NumberFormat = "@"
$cell.Value = "6.9%"
$cell.Style = "Normal"
$cell = $ws.Range("P12")
$cell.NumberFormat = "@"
$cell.Value = "30.9%"
$cell.Style = "Normal"
$cell = $ws.Range("P13")
$cell.NumberFormat = "@"
$cell.Value = "25.0%"
$cell.Style = "Normal"
$cell = $ws.Range("P14")
$cell.NumberFormat = "@"
$cell.Value = "11.6%"
$cell.Style = "Normal"
$cell = $ws.Range("P15")
$cell.NumberFormat = "@"
$cell.Value = "32.7%"
$cell.Style = "Normal"
$cell = $ws.Range("P16")
$cell.NumberFormat = "@"
$cell.Value = "30.7%"
$cell.Style = "Normal"
$cell = $ws.Range("P17")
$cell.NumberFormat = "@"
$cell.Value = "5.9%"
$cell.Style = "Normal"
$cell = $ws.Range("P18")
$cell.NumberFormat = "@"
$cell.Value = "27.9%"
$cell.Style = "Normal"
$cell = $ws.Range("P19")
$cell.NumberFormat = "@"
$cell.Value = "14.8%"
$cell.Style = "Normal"
$cell = $ws.Range("P20")
$cell.NumberFormat = "@"
$cell.Value = "3.7%"
$cell.Style = "Normal"
$cell = $ws.Range("P21")
$cell.NumberFormat = "@"
$cell.Value = "23.0%"
$cell.Style = "Normal"
$cell = $ws.Range("P22")
$cell.NumberFormat = "@"
$cell.Value = "11.7%"
$cell.Style = "Normal"

# Column Q, rows 2-22
$cell = $ws.Range("Q2")
$cell.NumberFormat = "@"
$cell.Value = "57.6%"
$cell.Style = "Normal"
$cell = $ws.Range("Q3")
$cell.NumberFormat = "@"
$cell.Value = "36.5%"
$cell.Style = "Normal"
$cell = $ws.Range("Q4")
$cell.NumberFormat = "@"
$cell.Value = "81.4%"
$cell.Style = "Normal"
$cell = $ws.Range("Q5")
$cell.NumberFormat = "@"
$cell.Value = "36.0%"
$cell.Style = "Normal"
$cell = $ws.Range("Q6")
$cell.NumberFormat = "@"
$cell.Value = "57.7%"
$cell.Style = "Normal"
$cell = $ws.Range("Q7")
$cell.NumberFormat = "@"
$cell.Value = "43.9%"
$cell.Style = "Normal"
$cell = $ws.Range("Q8")
$cell.NumberFormat = "@"
$cell.Value = "48.9%"
$cell.Style = "Normal"
$cell = $ws.Range("Q9")
$cell.NumberFormat = "@"
$cell.Value = "35.1%"
$cell.Style = "Normal"
$cell = $ws.Range("Q10")
$cell.NumberFormat = "@"
$cell.Value = "80.0%"
$cell.Style = "Normal"
$cell = $ws.Range("Q11")
$cell.NumberFormat = "@"
$cell.Value = "72.5%"
$cell.Style = "Normal"
$cell = $ws.Range("Q12")
$cell.NumberFormat = "@"
$cell.Value = "41.8%"
$cell.Style = "Normal"
$cell = $ws.Range("Q13")
$cell.NumberFormat = "@"
$cell.Value = "84.6%"
$cell.Style = "Normal"
$cell = $ws.Range("Q14")
$cell.NumberFormat = "@"
$cell.Value = "71.1%"
$cell.Style = "Normal"
$cell = $ws.Range("Q15")
$cell.NumberFormat = "@"
$cell.Value = "36.3%"
$cell.Style = "Normal"
$cell = $ws.Range("Q16")
$cell.NumberFormat = "@"
$cell.Value = "83.7%"
$cell.Style = "Normal"
$cell = $ws.Range("Q17")
$cell.NumberFormat = "@"
$cell.Value = "49.6%"
$cell.Style = "Normal"
$cell = $ws.Range("Q18")
$cell.NumberFormat = "@"
$cell.Value = "45.6%"
$cell.Style = "Normal"
$cell = $ws.Range("Q19")
$cell.NumberFormat = "@"
$cell.Value = "84.6%"
$cell.Style = "Normal"
$cell = $ws.Range("Q20")
$cell.NumberFormat = "@"
$cell.Value = "68.4%"
$cell.Style = "Normal"
$cell = $ws.Range("Q21")
$cell.NumberFormat = "@"
$cell.Value = "54.3%"
$cell.Style = "Normal"
$cell = $ws.Range("Q22")
$cell.NumberFormat = "@"
$cell.Value = "88.2%"
$cell.Style = "Normal"

